# Updates the cryptos list (rows 2-51) to reflect the latest scraped
# prices/volume percentages, and swaps the ARBITRUM / LidoDAOToken rows
# (35 and 36) to match the refreshed coinranking.com ordering.
#
# The "Price" column (D) holds values that look numeric (e.g. "0.9997",
# "241.49", "0.3110") but must stay stored as literal text - Excel's
# automatic type inference would otherwise coerce them to numbers and
# silently drop meaningful trailing zeros (e.g. "0.3110" -> 0.311) or
# collapse the thousand-dot-separated values (e.g. "29.213.66" would
# misparse). Forcing NumberFormat = "@" (Text) before the assignment
# keeps the literal string intact, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.213.66"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.42"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7079"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.49"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07811"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3110"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.85"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07808"
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.37"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.126"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.65"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6898"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.557"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008457"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.201.90"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "250.59"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.098.73"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.557"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1538"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.14"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.890"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.563"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.280"
$ws.Range("E30").Value = "  -2.87%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05213"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7603"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.853"
$ws.Range("E36").Value = "  -3.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.730"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8981"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.68"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.682"
$ws.Range("E44").Value = "  -11.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.997.87"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.26"
$ws.Range("E47").Value = "  -11.24%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.547"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.039"
$ws.Range("E51").Value = "  -0.69%  "
